# Update countries & provincias Spain
# - Refresh the "last updated" timestamp in A1
# - Two country rows swap rank order (new totals overtook the row above):
#     * Emiratos Arabes Unidos  <->  Paises Bajos     (rows 45/46)
#     * Corea del Sur           <->  Bosnia y Herzegovina (rows 76/77)
# - Refresh numeric stats (Casos totales / Nuevos casos / Casos activos /
#   Recuperados / Casos criticos / Muertes hoy / Muertes) for the rows whose
#   counters moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- "Datos actualizados ..." timestamp (row 1) ---------------------------
$ws.Range("A1").Value = "Datos actualizados a 2 de Septiembre de 2020 a las 15:38"

# --- Country label swaps ----------------------------------------------------
$ws.Range("A45").Value = "Paises Bajos"
$ws.Range("A46").Value = "Emiratos Arabes Unidos"

$ws.Range("A76").Value = "Bosnia y Herzegovina"
$ws.Range("A77").Value = "Corea del Sur"

# --- Numeric refresh --------------------------------------------------------
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6259973
$ws.Range("C4").Value = 2402
$ws.Range("D4").Value = 3497840
$ws.Range("E4").Value = 2573173
$ws.Range("G4").Value = 60
$ws.Range("H4").Value = 188960

# Row 6 - India
$ws.Range("B6").Value = 3794314
$ws.Range("C6").Value = 28206
$ws.Range("D6").Value = 2920122
$ws.Range("E6").Value = 807514
$ws.Range("G6").Value = 218
$ws.Range("H6").Value = 66678

# Row 13 - Argentina
$ws.Range("D13").Value = 315530
$ws.Range("E13").Value = 103738
$ws.Range("G13").Value = 52
$ws.Range("H13").Value = 8971

# Row 23 - Alemania
$ws.Range("B23").Value = 246234
$ws.Range("C23").Value = 233
$ws.Range("E23").Value = 15051
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 9383

# Row 42 - Suecia
$ws.Range("B42").Value = 84532
$ws.Range("G42").Value = 5
$ws.Range("H42").Value = 5820

# Row 44 - Bielorrusia
$ws.Range("B44").Value = 72141
$ws.Range("C44").Value = 179
$ws.Range("D44").Value = 70900
$ws.Range("E44").Value = 550
$ws.Range("G44").Value = 5
$ws.Range("H44").Value = 691

# Row 45 - Paises Bajos (after swap)
$ws.Range("B45").Value = 71863
$ws.Range("C45").Value = 734
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("G45").Value = 5
$ws.Range("H45").Value = 6235

# Row 46 - Emiratos Arabes Unidos (after swap)
$ws.Range("B46").Value = 71540
$ws.Range("C46").Value = 735
$ws.Range("D46").Value = 62029
$ws.Range("E46").Value = 9124
$ws.Range("G46").Value = 3
$ws.Range("H46").Value = 387

# Row 51 - Portugal
$ws.Range("B51").Value = 58633
$ws.Range("C51").Value = 390
$ws.Range("D51").Value = 42233
$ws.Range("E51").Value = 14573
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = 1827

# Row 67 - Azerbaiyan
$ws.Range("B67").Value = 36732
$ws.Range("C67").Value = 154
$ws.Range("D67").Value = 34116
$ws.Range("E67").Value = 2078
$ws.Range("G67").Value = 2
$ws.Range("H67").Value = 538

# Row 76 - Bosnia y Herzegovina (after swap)
$ws.Range("B76").Value = 20517
$ws.Range("C76").Value = 283
$ws.Range("D76").Value = 13626
$ws.Range("E76").Value = 6264
$ws.Range("G76").Value = 7
$ws.Range("H76").Value = 627

# Row 77 - Corea del Sur (after swap)
$ws.Range("B77").Value = 20449
$ws.Range("C77").Value = 267
$ws.Range("D77").Value = 15356
$ws.Range("E77").Value = 4767
$ws.Range("G77").Value = 2
$ws.Range("H77").Value = 326

# Row 86 - Republica de Macedonia
$ws.Range("B86").Value = 14600
$ws.Range("C86").Value = 145
$ws.Range("D86").Value = 11741
$ws.Range("E86").Value = 2253
$ws.Range("G86").Value = 2
$ws.Range("H86").Value = 606

# Row 164 - Vietnam
$ws.Range("B164").Value = 1046
$ws.Range("C164").Value = 2
$ws.Range("D164").Value = 746
$ws.Range("E164").Value = 266

# Row 184 - Gibraltar
$ws.Range("B184").Value = 295
$ws.Range("C184").Value = 5
$ws.Range("D184").Value = 240
$ws.Range("E184").Value = 55
